$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new cell values, row by row (A,B columns are unchanged; we are
# adding new "Firearms" (C) and "Sidearms" (D) columns). Filling in this
# row-major, left-to-right order keeps new shared-string entries appended
# in the same order Excel would naturally create them.
$ws.Range("C1").Value = "Firearms"
$ws.Range("D1").Value = "Sidearms"

$ws.Range("C2").Value = "Bad News"
$ws.Range("D2").Value = "Blunderbuss"

$ws.Range("C3").Value = "Blunderbuss"
$ws.Range("D3").Value = "Hand-Mortar"

$ws.Range("C4").Value = "Hand Mortar"
$ws.Range("D4").Value = "Palm Pistol"

$ws.Range("C5").Value = "Musket"
$ws.Range("D5").Value = "Pepperbox"

$ws.Range("C6").Value = "Palm Pistol"
$ws.Range("D6").Value = "Pistol"

$ws.Range("C7").Value = "Pepperbox"

$ws.Range("C8").Value = "Pistol"

# Copy the formatting from the existing "B" column cells (which already
# carry the sheet's shared cell style) onto the newly written cells so
# they pick up the same style instead of creating brand-new style records.
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("C3:D3").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("C4:D4").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("C5:D5").PasteSpecial(-4122)

$ws.Range("B6").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$ws.Range("B8").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$excel.CutCopyMode = 0
